# Add a "users" column to the "project hours" sheet, listing the users
# associated with each project.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# Header cell: same bold/bordered style as the other header cells (D1).
$ws.Range("E1").Value = "users"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: one user list per project, in the same row order as columns B-D.
$ws.Range("E2").Value = "['Hunter Young']"
$ws.Range("E3").Value = "['Arun Lakshmanan', 'Gabriel Barsi Haberfeld']"
$ws.Range("E4").Value = "['Man-Ki Yoon', 'Bo Liu']"
